# UPDATED ALL SS COURSES
# SS courses are now lec lab separated with new numbers from excel sheet SS COURSES.xlsx

$wb = $excel.ActiveWorkbook

$wsSS = $wb.Worksheets.Item("1171")   # SS courses sheet ("1171") - becomes the active sheet

# --- Row 15: catalog nbr 1006 -> 1007, LAB -> LEC ---
$wsSS.Range("F15").Value = 1007
$wsSS.Range("I15").Value = "LEC"

# --- Rows 71-78: catalog nbr 1008 -> 1009, LAB -> LEC ---
foreach ($r in 71..78) {
    $wsSS.Range("F$r").Value = 1009
    $wsSS.Range("I$r").Value = "LEC"
}

# --- Row 79: catalog nbr 1008 -> 1009, LAB -> LEC (also loses the bottom border -> style 7) ---
$wsSS.Range("F79").Value = 1009
$wsSS.Range("I79").Value = "LEC"
$wsSS.Range("I79").Borders.LineStyle = -4142   # xlLineStyleNone - drop bottom border like cols B:D,G,H,J,K,N keep (only I79 changes to no-border style)

# --- Rows 80-87: catalog nbr 1063 -> 1033 (component already LEC) ---
foreach ($r in 80..87) {
    $wsSS.Range("F$r").Value = 1033
}

# --- Rows 88-92: catalog nbr 1063 -> 1034, LAB -> LEC ---
foreach ($r in 88..92) {
    $wsSS.Range("F$r").Value = 1034
    $wsSS.Range("I$r").Value = "LEC"
}

# --- Sheet view / selection bookkeeping: "1171" becomes the active/selected sheet ---
# (its tabSelected flag turns on, "1168" loses it; the last touched cell is K87
#  and the view is scrolled down near the bottom of the SS block, topLeftCell A73)
$wsSS.Activate()
$excel.ActiveWindow.ScrollRow = 73
$wsSS.Range("K87").Select()
